# production shares with biomass shares
# Update existing CO2 utilization / bio-based feedstock rows with new values,
# and append four new "bio-based" combination rows below the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (rows 8-10) ---
$ws.Range("C8").Value = 3643380.558831645
$ws.Range("D8").Value = 5807844.520829624

$ws.Range("B9").Value = 8586549.999999989
$ws.Range("C9").Value = 957437.243744852
$ws.Range("D9").Value = 2249164.886954993

$ws.Range("C10").Value = 1680220.114095789
$ws.Range("D10").Value = 986796.9397839492

# --- Append new rows 13-16 ---
# Copy the formatting of the last labeled row (A12) onto the new label cells
# so the new rows reuse the existing bold/border/centered style.
$ws.Range("A12").Copy() | Out-Null

$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "Electrification + Bio-based feedstock"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0

$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "Conventional + Bio-based feedstock"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0

$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = "Conventional + Bio-based feedstock with CC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0

$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "Electrification + Bio-based"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0

$excel.CutCopyMode = $false
